# Update column F (dSF) values on Sheet1 to re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -4
    10 = -1
    14 = 2
    15 = 0
    24 = 5
    26 = 3
    27 = 3
    30 = 5
    31 = -1
    33 = -5
    39 = -4
    43 = 3
    44 = -7
    49 = -1
    51 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
